# Generate Report for Archive
# The localization CI regenerated the status report: two files that were
# previously "Ready for handoff" have now moved into translation, so their
# Status changes to "In Translation" on every sheet that tracks them. The
# dca1bbfd file's status is re-affirmed as "Ready for handoff" (unchanged).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
